$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1888586956521739
$ws.Range("C2").Value = 0.5665760869565217
$ws.Range("J2").Value = 0.02581521739130435
$ws.Range("O2").Value = 0.001358695652173913
$ws.Range("P2").Value = 0.1290760869565217
$ws.Range("S2").Value = 0.08831521739130435
$ws.Range("B3").Value = 0.007142857142857143
$ws.Range("C3").Value = 0.02142857142857143
$ws.Range("J3").Value = 0.02857142857142857
$ws.Range("P3").Value = 0.7428571428571429
$ws.Range("S3").Value = 0.2
$ws.Range("J4").Value = 0.04901960784313725
$ws.Range("P4").Value = 0.6666666666666666
$ws.Range("S4").Value = 0.2843137254901961
$ws.Range("P5").Value = 0.8333333333333334
$ws.Range("S5").Value = 0.1666666666666667
$ws.Range("B6").Value = 0.1051401869158878
$ws.Range("D6").Value = 0.02102803738317757
$ws.Range("F6").Value = 0.05841121495327103
$ws.Range("J6").Value = 0.2476635514018692
$ws.Range("O6").Value = 0.01869158878504673
$ws.Range("Q6").Value = 0.1425233644859813
$ws.Range("R6").Value = 0.05373831775700934
$ws.Range("S6").Value = 0.352803738317757
$ws.Range("B7").Value = 0.1381957773512476
$ws.Range("D7").Value = 0.02687140115163148
$ws.Range("E7").Value = 0.003838771593090211
$ws.Range("F7").Value = 0.04606525911708254
$ws.Range("J7").Value = 0.1285988483685221
$ws.Range("O7").Value = 0.01151631477927063
$ws.Range("Q7").Value = 0.1900191938579655
$ws.Range("R7").Value = 0.08445297504798464
$ws.Range("S7").Value = 0.3704414587332054
$ws.Range("B8").Value = 0.119960668633235
$ws.Range("D8").Value = 0.01966568338249754
$ws.Range("F8").Value = 0.04621435594886922
$ws.Range("J8").Value = 0.1170108161258604
$ws.Range("O8").Value = 0.01081612586037365
$ws.Range("Q8").Value = 0.1553588987217306
$ws.Range("R8").Value = 0.09734513274336283
$ws.Range("S8").Value = 0.4336283185840708
$ws.Range("B9").Value = 0.1088825214899714
$ws.Range("D9").Value = 0.02865329512893983
$ws.Range("F9").Value = 0.07736389684813753
$ws.Range("J9").Value = 0.1174785100286533
$ws.Range("O9").Value = 0.01146131805157593
$ws.Range("Q9").Value = 0.1805157593123209
$ws.Range("R9").Value = 0.09742120343839542
$ws.Range("S9").Value = 0.3782234957020058
$ws.Range("B10").Value = 0.1251509054325956
$ws.Range("D10").Value = 0.02293762575452716
$ws.Range("E10").Value = 0.001609657947686117
$ws.Range("F10").Value = 0.06841046277665996
$ws.Range("J10").Value = 0.1138832997987928
$ws.Range("O10").Value = 0.01569416498993964
$ws.Range("Q10").Value = 0.2056338028169014
$ws.Range("R10").Value = 0.09255533199195171
$ws.Range("S10").Value = 0.3541247484909457
$ws.Range("G11").Value = 0.1516034985422741
$ws.Range("J11").Value = 0.08017492711370262
$ws.Range("K11").Value = 0.1909620991253644
$ws.Range("L11").Value = 0.5626822157434402
$ws.Range("S11").Value = 0.01457725947521866
$ws.Range("G12").Value = 0.7971014492753623
$ws.Range("J12").Value = 0.1570048309178744
$ws.Range("K12").Value = 0.004830917874396135
$ws.Range("L12").Value = 0.01932367149758454
$ws.Range("S12").Value = 0.02173913043478261
$ws.Range("G13").Value = 0.725
$ws.Range("J13").Value = 0.2166666666666667
$ws.Range("S13").Value = 0.05833333333333333
$ws.Range("F15").Value = 0.024330900243309
$ws.Range("H15").Value = 0.218978102189781
$ws.Range("I15").Value = 0.0583941605839416
$ws.Range("J15").Value = 0.3114355231143552
$ws.Range("K15").Value = 0.06326034063260341
$ws.Range("M15").Value = 0.0194647201946472
$ws.Range("N15").Value = 0.0024330900243309
$ws.Range("O15").Value = 0.072992700729927
$ws.Range("S15").Value = 0.2287104622871046
$ws.Range("F16").Value = 0.02164502164502164
$ws.Range("H16").Value = 0.1731601731601732
$ws.Range("I16").Value = 0.08441558441558442
$ws.Range("J16").Value = 0.3506493506493507
$ws.Range("K16").Value = 0.1233766233766234
$ws.Range("M16").Value = 0.0303030303030303
$ws.Range("N16").Value = 0.006493506493506494
$ws.Range("O16").Value = 0.0670995670995671
$ws.Range("S16").Value = 0.1428571428571428
$ws.Range("F17").Value = 0.02015677491601344
$ws.Range("H17").Value = 0.2138857782754759
$ws.Range("I17").Value = 0.05487122060470324
$ws.Range("J17").Value = 0.4053751399776036
$ws.Range("K17").Value = 0.1052631578947368
$ws.Range("M17").Value = 0.02799552071668533
$ws.Range("N17").Value = 0.001119820828667413
$ws.Range("O17").Value = 0.05375139977603583
$ws.Range("S17").Value = 0.1175811870100784
$ws.Range("F18").Value = 0.0162037037037037
$ws.Range("H18").Value = 0.1944444444444444
$ws.Range("I18").Value = 0.07175925925925926
$ws.Range("J18").Value = 0.3981481481481481
$ws.Range("K18").Value = 0.1157407407407407
$ws.Range("M18").Value = 0.02546296296296296
$ws.Range("N18").Value = 0.002314814814814815
$ws.Range("O18").Value = 0.03703703703703703
$ws.Range("S18").Value = 0.1388888888888889
$ws.Range("F19").Value = 0.01117964533538936
$ws.Range("H19").Value = 0.2197378565921357
$ws.Range("I19").Value = 0.07787201233616037
$ws.Range("J19").Value = 0.3492675404780262
$ws.Range("K19").Value = 0.1214340786430224
$ws.Range("M19").Value = 0.02467232074016962
$ws.Range("N19").Value = 0.001156515034695451
$ws.Range("O19").Value = 0.0670995670995671
$ws.Range("S19").Value = 0.1283731688511951
